$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (currently right after the
#    "N 2" run in the first heading paragraph).
# ---------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------
# 2. Pin-table text shifts (Table 2, column 1): each label moves down
#    one row, "18" (with the lastRenderedPageBreak run) becomes "RXD".
# ---------------------------------------------------------------------
$tbl = $d.Tables.Item(2)
$tbl.Cell(6, 1).Range.Text  = "RXD"
$tbl.Cell(7, 1).Range.Text  = "18"
$tbl.Cell(8, 1).Range.Text  = "27"
$tbl.Cell(9, 1).Range.Text  = "24"
$tbl.Cell(10, 1).Range.Text = "MOSI"

# ---------------------------------------------------------------------
# 3. Re-insert "_GoBack" right after "D15_6" (last row of Table 2,
#    2nd column), i.e. at the very end of that cell's text, before the
#    paragraph mark.
#
#    A collapsed Range sitting exactly on a paragraph-end boundary
#    confuses this runtime's Bookmarks.Add (it silently resets to the
#    start of the document), so we dodge the edge case: temporarily
#    insert a placeholder character after the target text, add the
#    bookmark just in front of it (a perfectly safe, non-boundary
#    position), then delete the placeholder again. The bookmark stays
#    anchored where we put it.
# ---------------------------------------------------------------------
$lastCell = $tbl.Cell($tbl.Rows.Count, 2)
$cellRange = $lastCell.Range
$textEnd = $cellRange.End - 1

$placeholder = $d.Range($textEnd, $textEnd)
$placeholder.InsertAfter("X")

$bmPos = $d.Range($textEnd, $textEnd)
$d.Bookmarks.Add("_GoBack", $bmPos)

$placeholderRange = $d.Range($textEnd, $textEnd + 1)
$placeholderRange.Delete()
